$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.876.79'
$ws.Range('E2').Value = '  -2.31%  '
$ws.Range('D3').Value = '1.784.40'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('E4').Value = '  +0.80%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.49'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.009'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3606'
$ws.Range('E8').Value = '  -1.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07146'
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8391'
$ws.Range('E10').Value = '  -3.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.23'
$ws.Range('E11').Value = '  -2.11%  '
$ws.Range('D12').Value = '1.834.22'
$ws.Range('E12').Value = '  -5.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.244'
$ws.Range('E13').Value = '  -3.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.330'
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06808'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.28'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.91'
$ws.Range('E20').Value = '  -3.31%  '
$ws.Range('D21').Value = '27.019.20'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.031'
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.03'
$ws.Range('E23').Value = '  +1.90%  '
$ws.Range('D24').Value = '2.033.17'
$ws.Range('E24').Value = '  -4.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.928'
$ws.Range('E25').Value = '  -2.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.08'
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.14'
$ws.Range('E27').Value = '  -3.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.06'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.003'
$ws.Range('E29').Value = '  -2.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.627'
$ws.Range('E30').Value = '  -11.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08948'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7135'
$ws.Range('E32').Value = '  -5.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.847'
$ws.Range('E33').Value = '  -4.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.304'
$ws.Range('E34').Value = '  -5.35%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.009'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.084'
$ws.Range('E36').Value = '  -4.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.077'
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('E39').Value = '  -4.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4927'
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1611'
$ws.Range('E41').Value = '  -3.47%  '
$ws.Range('E42').Value = '  -10.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.979'
$ws.Range('E43').Value = '  -9.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.906'
$ws.Range('E44').Value = '  -5.84%  '
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.10'
$ws.Range('E46').Value = '  -1.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.05'
$ws.Range('E47').Value = '  -3.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06259'
$ws.Range('E48').Value = '  -3.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4468'
$ws.Range('E49').Value = '  -4.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.572'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.702'
$ws.Range('E51').Value = '  -2.05%  '
